$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B ("Nombre proveedor"), shifting
# the existing B:E columns to C:F.
$ws.Columns("B:B").Insert()

# Give the new header cell (B1) the same formatting as the other
# header cells (bold + border), then set its text.
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B1").Value = "producto"

# Give the new data cells (B2:B29) the same plain formatting as the
# neighbouring text column (C2:C29) instead of the inherited date
# format, then fill them with the new "producto" value.
$ws.Range("C2:C29").Copy()
$ws.Range("B2:B29").PasteSpecial(-4122)  # xlPasteFormats
for ($r = 2; $r -le 29; $r++) {
    $ws.Cells.Item($r, 2).Value = "DESCONOCIDO"
}

$excel.CutCopyMode = 0
